$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.881601
$ws.Range("H2").Value = 2.644803
$ws.Range("I2").Value = 0.02072192623875441
$ws.Range("J2").Value = 0.02072192623875442
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.976167666666666
$ws.Range("N2").Value = 14.928503
$ws.Range("O2").Value = 0.195706700443638
$ws.Range("P2").Value = 0.195706700443638
$ws.Range("Q2").Value = 4.386994391100999
$ws.Range("R2").Value = 39.482949519909
$ws.Range("S2").Value = 0.004055419811023073
$ws.Range("T2").Value = 0.004055419811023074

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.881601
$ws.Range("H3").Value = 2.644803
$ws.Range("I3").Value = 0.02072192623875441
$ws.Range("J3").Value = 0.02072192623875442
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.170105
$ws.Range("N3").Value = 21.510315
$ws.Range("O3").Value = 0.2819916219431576
$ws.Range("P3").Value = 0.2819916219431576
$ws.Range("Q3").Value = 6.321171738104999
$ws.Range("R3").Value = 56.890545642945
$ws.Range("S3").Value = 0.005843409589852832
$ws.Range("T3").Value = 0.005843409589852833

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.881601
$ws.Range("H4").Value = 2.644803
$ws.Range("I4").Value = 0.02072192623875441
$ws.Range("J4").Value = 0.02072192623875442
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.59057
$ws.Range("N4").Value = 4.77171
$ws.Range("O4").Value = 0.06255520862164893
$ws.Range("P4").Value = 0.06255520862164894
$ws.Range("Q4").Value = 1.40224810257
$ws.Range("R4").Value = 12.62023292313
$ws.Range("S4").Value = 0.001296264418907703
$ws.Range("T4").Value = 0.001296264418907704

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.881601
$ws.Range("H5").Value = 2.644803
$ws.Range("I5").Value = 0.02072192623875441
$ws.Range("J5").Value = 0.02072192623875442
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.818291
$ws.Range("N5").Value = 8.454873
$ws.Range("O5").Value = 0.1108400016733093
$ws.Range("P5").Value = 0.1108400016733093
$ws.Range("Q5").Value = 2.484608163891
$ws.Range("R5").Value = 22.361473475019
$ws.Range("S5").Value = 0.002296818338977732
$ws.Range("T5").Value = 0.002296818338977732

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.881601
$ws.Range("H6").Value = 2.644803
$ws.Range("I6").Value = 0.02072192623875441
$ws.Range("J6").Value = 0.02072192623875442
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.871526000000001
$ws.Range("N6").Value = 26.614578
$ws.Range("O6").Value = 0.3489064673182461
$ws.Range("P6").Value = 0.3489064673182461
$ws.Range("Q6").Value = 7.821146193126001
$ws.Range("R6").Value = 70.390315738134
$ws.Range("S6").Value = 0.007230014079993075
$ws.Range("T6").Value = 0.007230014079993075

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.27504099999999
$ws.Range("H7").Value = 111.825123
$ws.Range("I7").Value = 0.8761453879346173
$ws.Range("J7").Value = 0.8761453879346174
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.976167666666666
$ws.Range("N7").Value = 14.928503
$ws.Range("O7").Value = 0.195706700443638
$ws.Range("P7").Value = 0.195706700443638
$ws.Range("Q7").Value = 185.4868537978743
$ws.Range("R7").Value = 1669.381684180869
$ws.Range("S7").Value = 0.1714675229815952
$ws.Range("T7").Value = 0.1714675229815952

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 37.27504099999999
$ws.Range("H8").Value = 111.825123
$ws.Range("I8").Value = 0.8761453879346173
$ws.Range("J8").Value = 0.8761453879346174
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.170105
$ws.Range("N8").Value = 21.510315
$ws.Range("O8").Value = 0.2819916219431576
$ws.Range("P8").Value = 0.2819916219431576
$ws.Range("Q8").Value = 267.2659578493049
$ws.Range("R8").Value = 2405.393620643745
$ws.Range("S8").Value = 0.2470656590016997
$ws.Range("T8").Value = 0.2470656590016998

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 37.27504099999999
$ws.Range("H9").Value = 111.825123
$ws.Range("I9").Value = 0.8761453879346173
$ws.Range("J9").Value = 0.8761453879346174
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.59057
$ws.Range("N9").Value = 4.77171
$ws.Range("O9").Value = 0.06255520862164893
$ws.Range("P9").Value = 0.06255520862164894
$ws.Range("Q9").Value = 59.28856196336999
$ws.Range("R9").Value = 533.59705767033
$ws.Range("S9").Value = 0.05480745752514551
$ws.Range("T9").Value = 0.05480745752514554

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.27504099999999
$ws.Range("H10").Value = 111.825123
$ws.Range("I10").Value = 0.8761453879346173
$ws.Range("J10").Value = 0.8761453879346174
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.818291
$ws.Range("N10").Value = 8.454873
$ws.Range("O10").Value = 0.1108400016733093
$ws.Range("P10").Value = 0.1108400016733093
$ws.Range("Q10").Value = 105.051912574931
$ws.Range("R10").Value = 945.4672131743788
$ws.Range("S10").Value = 0.09711195626473522
$ws.Range("T10").Value = 0.09711195626473523

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 37.27504099999999
$ws.Range("H11").Value = 111.825123
$ws.Range("I11").Value = 0.8761453879346173
$ws.Range("J11").Value = 0.8761453879346174
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.871526000000001
$ws.Range("N11").Value = 26.614578
$ws.Range("O11").Value = 0.3489064673182461
$ws.Range("P11").Value = 0.3489064673182461
$ws.Range("Q11").Value = 330.686495382566
$ws.Range("R11").Value = 2976.178458443094
$ws.Range("S11").Value = 0.3056927921614416
$ws.Range("T11").Value = 0.3056927921614417

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.387713666666667
$ws.Range("H12").Value = 13.163141
$ws.Range("I12").Value = 0.1031326858266283
$ws.Range("J12").Value = 0.1031326858266283
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.976167666666666
$ws.Range("N12").Value = 14.928503
$ws.Range("O12").Value = 0.195706700443638
$ws.Range("P12").Value = 0.195706700443638
$ws.Range("Q12").Value = 21.83399887865811
$ws.Range("R12").Value = 196.505989907923
$ws.Range("S12").Value = 0.02018375765101978
$ws.Range("T12").Value = 0.02018375765101978

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.387713666666667
$ws.Range("H13").Value = 13.163141
$ws.Range("I13").Value = 0.1031326858266283
$ws.Range("J13").Value = 0.1031326858266283
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.170105
$ws.Range("N13").Value = 21.510315
$ws.Range("O13").Value = 0.2819916219431576
$ws.Range("P13").Value = 0.2819916219431576
$ws.Range("Q13").Value = 31.460367699935
$ws.Range("R13").Value = 283.143309299415
$ws.Range("S13").Value = 0.02908255335160502
$ws.Range("T13").Value = 0.02908255335160502

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.387713666666667
$ws.Range("H14").Value = 13.163141
$ws.Range("I14").Value = 0.1031326858266283
$ws.Range("J14").Value = 0.1031326858266283
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.59057
$ws.Range("N14").Value = 4.77171
$ws.Range("O14").Value = 0.06255520862164893
$ws.Range("P14").Value = 0.06255520862164894
$ws.Range("Q14").Value = 6.978965726789999
$ws.Range("R14").Value = 62.81069154111
$ws.Range("S14").Value = 0.00645148667759571
$ws.Range("T14").Value = 0.006451486677595711

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.387713666666667
$ws.Range("H15").Value = 13.163141
$ws.Range("I15").Value = 0.1031326858266283
$ws.Range("J15").Value = 0.1031326858266283
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.818291
$ws.Range("N15").Value = 8.454873
$ws.Range("O15").Value = 0.1108400016733093
$ws.Range("P15").Value = 0.1108400016733093
$ws.Range("Q15").Value = 12.36585393734367
$ws.Range("R15").Value = 111.292685436093
$ws.Range("S15").Value = 0.01143122706959637
$ws.Range("T15").Value = 0.01143122706959637

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.387713666666667
$ws.Range("H16").Value = 13.163141
$ws.Range("I16").Value = 0.1031326858266283
$ws.Range("J16").Value = 0.1031326858266283
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 8.871526000000001
$ws.Range("N16").Value = 26.614578
$ws.Range("O16").Value = 0.3489064673182461
$ws.Range("P16").Value = 0.3489064673182461
$ws.Range("Q16").Value = 38.92571587438867
$ws.Range("R16").Value = 350.331442869498
$ws.Range("S16").Value = 0.03598366107681144
$ws.Range("T16").Value = 0.03598366107681144
